# Append 3 blank paragraphs and a 4th paragraph containing the heatmap
# code snippet after the document's final paragraph
# ("transitions.to_csv(...)").

$d = $word.ActiveDocument

$range = $d.Paragraphs.Last.Range
$range.Collapse(0)
$range.InsertParagraphAfter()
$range.Collapse(0)
$range.InsertParagraphAfter()
$range.Collapse(0)
$range.InsertParagraphAfter()
$range.Collapse(0)
$range.InsertParagraphAfter()

$heatmapText = "# Step 1: Pivot Data for Heatmap heatmap_data = df.pivot(index='journey_name', columns='path', values='count') heatmap_data.fillna(0, inplace=True) # Replace NaNs with 0 for missing counts # Step 2: Plot Heatmap plt.figure(figsize=(12, 8)) sns.heatmap(heatmap_data, annot=True, fmt=`".0f`", cmap=`"YlGnBu`", cbar_kws={'label': 'Frequency'}) plt.title(`"Heatmap of Paths for Each Journey`") plt.xlabel(`"Paths`") plt.ylabel(`"Journeys`") plt.xticks(rotation=45, ha=`"right`") plt.tight_layout() plt.show() "

$d.Paragraphs.Last.Range.Text = $heatmapText

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
